# Updates cryptos list values (Price / Volume(1h) columns, and the
# Aave/Mantle row swap) to match the refreshed data snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '59.221.50'
$ws.Range("E2").Value = '  -0.36%  '
# Row 3
$ws.Range("D3").Value = '2.525.98'
$ws.Range("E3").Value = '  +0.16%  '
# Row 4
$ws.Range("E4").Value = '  +0.28%  '
# Row 5
$ws.Range("D5").Value = '''536.61'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.14%  '
# Row 6
$ws.Range("D6").Value = '''136.87'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.52%  '
# Row 7
$ws.Range("E7").Value = '  +0.19%  '
# Row 8
$ws.Range("D8").Value = '''0.570'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.80%  '
# Row 9
$ws.Range("D9").Value = '2.524.01'
$ws.Range("E9").Value = '  -0.17%  '
# Row 10
$ws.Range("E10").Value = '  -0.32%  '
# Row 11
$ws.Range("E11").Value = '  -2.47%  '
# Row 12
$ws.Range("D12").Value = '''5.30'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.44%  '
# Row 13
$ws.Range("E13").Value = '  -0.92%  '
# Row 14
$ws.Range("D14").Value = '2.974.26'
$ws.Range("E14").Value = '  +0.01%  '
# Row 15
$ws.Range("D15").Value = '''23.07'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.58%  '
# Row 16
$ws.Range("D16").Value = '59.145.04'
$ws.Range("E16").Value = '  -0.33%  '
# Row 17
$ws.Range("E17").Value = '  -1.73%  '
# Row 18
$ws.Range("D18").Value = '2.521.07'
$ws.Range("E18").Value = '  +0.08%  '
# Row 19
$ws.Range("D19").Value = '''11.14'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.22%  '
# Row 20
$ws.Range("E20").Value = '  -0.17%  '
# Row 21
$ws.Range("D21").Value = '''324.04'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.58%  '
# Row 22
$ws.Range("D22").Value = '''1.00'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.09%  '
# Row 23
$ws.Range("D23").Value = '''5.95'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.20%  '
# Row 24
$ws.Range("D24").Value = '''65.52'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.44%  '
# Row 25
$ws.Range("E25").Value = '  -0.24%  '
# Row 26
$ws.Range("E26").Value = '  -2.29%  '
# Row 27
$ws.Range("E27").Value = '  -0.05%  '
# Row 28
$ws.Range("E28").Value = '  -2.70%  '
# Row 29
$ws.Range("D29").Value = '0.0₃0773'
$ws.Range("E29").Value = '  -1.64%  '
# Row 30
$ws.Range("E30").Value = '  -0.72%  '
# Row 31
$ws.Range("E31").Value = '  -1.98%  '
# Row 32
$ws.Range("D32").Value = '''171.65'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.82%  '
# Row 33
$ws.Range("E33").Value = '  +5.64%  '
# Row 34
$ws.Range("D34").Value = '''0.999'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.05%  '
# Row 35
$ws.Range("E35").Value = '  -0.86%  '
# Row 36
$ws.Range("D36").Value = '''18.38'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.02%  '
# Row 37
$ws.Range("E37").Value = '  -2.29%  '
# Row 38
$ws.Range("D38").Value = '''1.54'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.26%  '
# Row 39
$ws.Range("D39").Value = '''36.71'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.78%  '
# Row 40
$ws.Range("D40").Value = '''0.813'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.29%  '
# Row 41
$ws.Range("D41").Value = '''3.60'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.20%  '
# Row 42
$ws.Range("D42").Value = '''285.38'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.26%  '
# Row 43
$ws.Range("D43").Value = '''5.11'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.39%  '
# Row 44
$ws.Range("E44").Value = '  +0.18%  '
# Row 45
$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").Value = '''0.610'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.35%  '
# Row 46
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '''131.21'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.12%  '
# Row 47
$ws.Range("D47").Value = '''10.88'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.17%  '
# Row 48
$ws.Range("E48").Value = '  -1.64%  '
# Row 49
$ws.Range("D49").Value = '''0.0507'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.30%  '
# Row 50
$ws.Range("E50").Value = '  -1.68%  '
# Row 51
$ws.Range("D51").Value = '''17.41'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.03%  '
